$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap 'sum' and 'author5' columns (S <-> T) ---
$ws.Range("S1").Value = 'author5'
$ws.Range("T1").Value = 'sum'

# --- Row 2: author moves to S, summary splits into paragraphs T:W ---
$ws.Range("S2").Value = 'Anonym'
$ws.Range("T2").Value = 'Die Nutzung der App ist aus der Sicht von Experten aus weiteren Gründen problematisch. Nicht nur, weil TikTok sehr viele Daten über seine Nutzer sammelt, sondern auch, weil eine übermäßige Nutzung der Plattform zu massivem Suchtverhalten führen kann, welches den gesamten Alltag von Jugendlichen beherrschen kann.'
$ws.Range("U2").Value = 'So wird beispielsweise die Kapazität deines Arbeitsgedächtnisses durch TikTok-Nutzern reduziert. Das bedeutet, dass man dadurch aktiv „verdummt“. Deshalb hat man oft keine Lust mehr, über einen längeren Zeitraum etwas anderes zu tun oder sich zu konzentrieren. Man kann/will z.B. Konversationen nicht mehr aufmerksam zuhören. Das alles funktioniert, weil die Nutzer durch TikTok kurzfristige positive Rückmeldungen bekommen, wie bei einer Art Belohnung. Dabei sind es im wirklichen Leben die langfristigen Belohnungen, die einen glücklich machen, z.B. mit der Familie etwas zu unternehmen. '
$ws.Range("V2").Value = 'Nicht nur die Art, wie TikTok seine User an sich bindet, kann schädlich sein, sondern auch, welcher Content auf den Nutzer zugeschnitten wird. Ein großer Nachteil ist zudem, dass bei süchtigen TikTok-Nutzern Depressionen, Angst und Stress auftreten können. Wenn man z.B. schon solche Tendenzen zeigt, werden Stimmungen wie Traurigkeit durch thematisch zugeschnittene Videos auf TikTok zusätzlich unterstützt. Damit ist TikTok nicht nur für die mentale Gesundheit des Einzelnen, sondern auch für uns als Gesellschaft eine große Gefahr.'
$ws.Range("W2").Value = 'Dass TikTok weiterreichende Auswirkungen auf unsere Jugend, unsere Psyche und unsere Gesellschaft hat, sollte dir spätestens jetzt bewusst sein. Was hältst du von TikTok? Schreib‘ uns deine Meinung: schuelerzeitung@gmg.amberg.de'

# --- Rows 3-6: author moves to S, summary text moves to T ---
$ws.Range("S3").Value = 'Anonym'
$ws.Range("T3").Value = 'summary 2'
$ws.Range("S4").Value = 'Anonym'
$ws.Range("T4").Value = 'summary 3'
$ws.Range("S5").Value = 'Anonym'
$ws.Range("T5").Value = 'summary 4'
$ws.Range("S6").Value = 'Anonym'
$ws.Range("T6").Value = 'summary 5'

# --- Rename the built-in cell style from German 'Standard' to 'Normal' ---
$wb.Styles.Item("Standard").Name = "Normal"

# --- Update selection to reflect new active cell ---
$ws.Range("W4").Select() | Out-Null
